# Auto-applied price/profit refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# (scheduled market-data runner update)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 572
$ws.Range("I2").Value = 650
$ws.Range("J2").Value = 299
$ws.Range("K2").Value = 650
$ws.Range("L2").Value = 299
$ws.Range("M2").Value = -537
$ws.Range("N2").Value = -525
$ws.Range("H33").Value = 679
$ws.Range("I33").Value = 205
$ws.Range("K33").Value = 205
$ws.Range("M33").Value = 24
$ws.Range("H40").Value = 3703.9473
$ws.Range("I40").Value = 2852.4546
$ws.Range("K40").Value = 2852.4546
$ws.Range("M40").Value = -2677.4546
$ws.Range("H125").Value = 3080.647
$ws.Range("I125").Value = 2337.1
$ws.Range("K125").Value = 21033.9
$ws.Range("M125").Value = -18573.9
$ws.Range("H138").Value = 4591.1265
$ws.Range("J138").Value = 5202.9253
$ws.Range("L138").Value = 15608.7759
$ws.Range("N138").Value = -25888.7759
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1229.5
$ws.Range("I26").Value = 1536.6666
$ws.Range("J26").Value = 308
$ws.Range("K26").Value = 1536.6666
$ws.Range("L26").Value = 308
$ws.Range("M26").Value = -1206.6666
$ws.Range("N26").Value = -968
$ws.Range("H43").Value = 22342
$ws.Range("I43").Value = 22342
$ws.Range("K43").Value = 22342
$ws.Range("M43").Value = -22029
$ws.Range("H61").Value = 4257.4
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents() | Out-Null
$ws.Range("H74").Value = 2534.8948
$ws.Range("I74").Value = 2116.077
$ws.Range("J74").Value = 3442.3333
$ws.Range("K74").Value = 2116.077
$ws.Range("L74").Value = 3442.3333
$ws.Range("M74").Value = -1242.077
$ws.Range("N74").Value = -5190.3333
$ws.Range("H77").Value = 2534.8948
$ws.Range("I77").Value = 2116.077
$ws.Range("J77").Value = 3442.3333
$ws.Range("K77").Value = 10580.385
$ws.Range("L77").Value = 17211.6665
$ws.Range("M77").Value = -6212.385000000002
$ws.Range("N77").Value = -25947.6665
$ws.Range("H102").Value = 3615.6428
$ws.Range("I102").Value = 3122
$ws.Range("K102").Value = 3122
$ws.Range("M102").Value = -1500
$ws.Range("H110").Value = 10070.917
$ws.Range("I110").Value = 10804.637
$ws.Range("K110").Value = 10804.637
$ws.Range("M110").Value = -8759.637000000001
$ws.Range("H136").Value = 4257.4
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents() | Out-Null
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7019.963
$ws.Range("I134").Value = 7097.654
$ws.Range("K134").Value = 21292.962
$ws.Range("M134").Value = -18757.962
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7712.1816
$ws.Range("I31").Value = 6639
$ws.Range("K31").Value = 6639
$ws.Range("M31").Value = -6344
$ws.Range("H34").Value = 7712.1816
$ws.Range("I34").Value = 6639
$ws.Range("K34").Value = 6639
$ws.Range("M34").Value = -6437
$ws.Range("H99").Value = 282751.28
$ws.Range("J99").Value = 5248.4443
$ws.Range("L99").Value = 5248.4443
$ws.Range("N99").Value = -8244.444299999999
$ws.Range("H126").Value = 282751.28
$ws.Range("J126").Value = 5248.4443
$ws.Range("L126").Value = 15745.3329
$ws.Range("N126").Value = -20685.3329
$ws.Range("H132").Value = 31364.611
$ws.Range("I132").Value = 9680.177
$ws.Range("J132").Value = 400000
$ws.Range("K132").Value = 29040.531
$ws.Range("L132").Value = 1200000
$ws.Range("M132").Value = -26510.531
$ws.Range("N132").Value = -1205060
$ws.Range("H134").Value = 1333755.5
$ws.Range("I134").Value = 1606842.4
$ws.Range("J134").Value = 2457.125
$ws.Range("K134").Value = 4820527.199999999
$ws.Range("L134").Value = 7371.375
$ws.Range("M134").Value = -4817992.199999999
$ws.Range("N134").Value = -12441.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 3634.7273
$ws.Range("I18").Value = 4268.75
$ws.Range("J18").Value = 1944
$ws.Range("K18").Value = 12806.25
$ws.Range("L18").Value = 5832
$ws.Range("M18").Value = -12637.25
$ws.Range("N18").Value = -6170
$ws.Range("H122").Value = 1151.8572
$ws.Range("I122").Value = 807.2308
$ws.Range("J122").Value = 1450.5333
$ws.Range("K122").Value = 7265.077200000001
$ws.Range("L122").Value = 13054.7997
$ws.Range("M122").Value = -4815.077200000001
$ws.Range("N122").Value = -17954.7997
$ws.Range("H129").Value = 27779140
$ws.Range("J129").Value = 55557736
$ws.Range("L129").Value = 166673208
$ws.Range("N129").Value = -166683208
$ws.Range("H132").Value = 40509.08
$ws.Range("I132").Value = 858.5
$ws.Range("J132").Value = 58131.555
$ws.Range("K132").Value = 7726.5
$ws.Range("L132").Value = 523183.995
$ws.Range("M132").Value = -5196.5
$ws.Range("N132").Value = -528243.995
$ws.Range("H134").Value = 3610
$ws.Range("I134").Value = 830
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 2490
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 2580
$ws.Range("N134").Value = -25140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9663.583000000001
$ws.Range("I70").Value = 9000.833000000001
$ws.Range("J70").Value = 10326.333
$ws.Range("K70").Value = 9000.833000000001
$ws.Range("L70").Value = 10326.333
$ws.Range("M70").Value = -8730.833000000001
$ws.Range("N70").Value = -10866.333
$ws.Range("H73").Value = 9663.583000000001
$ws.Range("I73").Value = 9000.833000000001
$ws.Range("J73").Value = 10326.333
$ws.Range("K73").Value = 9000.833000000001
$ws.Range("L73").Value = 10326.333
$ws.Range("M73").Value = -8064.833000000001
$ws.Range("N73").Value = -12198.333
$ws.Range("H80").Value = 3108.125
$ws.Range("I80").Value = 2994.75
$ws.Range("J80").Value = 3221.5
$ws.Range("K80").Value = 2994.75
$ws.Range("L80").Value = 3221.5
$ws.Range("M80").Value = -1996.75
$ws.Range("N80").Value = -5217.5
$ws.Range("H83").Value = 3108.125
$ws.Range("I83").Value = 2994.75
$ws.Range("J83").Value = 3221.5
$ws.Range("K83").Value = 14973.75
$ws.Range("L83").Value = 16107.5
$ws.Range("M83").Value = -9981.75
$ws.Range("N83").Value = -26091.5
$ws.Range("H102").Value = 9963.429
$ws.Range("I102").Value = 12959.8
$ws.Range("J102").Value = 2472.5
$ws.Range("K102").Value = 12959.8
$ws.Range("L102").Value = 2472.5
$ws.Range("M102").Value = -11337.8
$ws.Range("N102").Value = -5716.5
$ws.Range("H125").Value = 61108.332
$ws.Range("J125").Value = 61108.332
$ws.Range("L125").Value = 61108.332
$ws.Range("N125").Value = -66028.33199999999
$ws.Range("H126").Value = 24528.691
$ws.Range("I126").Value = 43742
$ws.Range("J126").Value = 15989.444
$ws.Range("K126").Value = 131226
$ws.Range("L126").Value = 47968.33199999999
$ws.Range("M126").Value = -128756
$ws.Range("N126").Value = -52908.33199999999
$ws.Range("H132").Value = 3777
$ws.Range("I132").Value = 2318.2
$ws.Range("J132").Value = 9247.5
$ws.Range("K132").Value = 6954.599999999999
$ws.Range("L132").Value = 27742.5
$ws.Range("M132").Value = -4424.599999999999
$ws.Range("N132").Value = -32802.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 371.83334
$ws.Range("I16").Value = 371.83334
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 371.83334
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -201.83334
$ws.Range("N16").ClearContents() | Out-Null
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2927.3076
$ws.Range("I107").Value = 4694.25
$ws.Range("K107").Value = 14082.75
$ws.Range("M107").Value = -12162.75
$ws.Range("H136").Value = 2088
$ws.Range("I136").Value = 1455.2632
$ws.Range("K136").Value = 4365.7896
$ws.Range("M136").Value = -1815.7896
